$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "#7cc867#fb5b89#c885da",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "#7cc867: 32^p#fb5b89: 13^p#c885da: 12",
    2
)
